# "New giant refactor and cleanup"
#
# Semantic changes applied to data/guns.xlsx (weapons.csv sheet):
#   - I1 header renamed:  "Shots"  -> "XShots"
#   - H2 ammo value fixed: "10mm"  -> "10-mm"
#   - Selection/active-cell cursor moved from L3 to D16
#   - Workbook window chrome (position/size) updated to new values
#     (best-effort; Excel view-state only, no data impact)
#
# Everything else in the underlying OOXML diff (re-ordering of the
# sharedStrings table, renumbered <v> indices, the dropped <ignoredErrors>
# block) is a pure side effect of Excel re-serialising the file and carries
# no data meaning - it is not something the object model exposes, so it is
# intentionally not reproduced here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cell content changes -------------------------------------------------
$ws.Range("I1").Value = "XShots"
$ws.Range("H2").Value = "10-mm"

# --- cursor / selection ----------------------------------------------------
$null = $ws.Range("D16").Select()

# --- workbook window chrome (position/size) --------------------------------
$win = $excel.ActiveWindow
$win.Left = 3920
$win.Top = 4220
$win.Width = 19000
$win.Height = 15160
